# Data cleaning: remove the duplicate customer case row (Alice Johnson /
# CS2025_00073) that was accidentally re-entered as a near-duplicate of the
# first record, then renumber the remaining CaseID values sequentially
# starting at CS2025_00142 for the (OpenAI response) training data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 holds the duplicate "Alice Johnson" / "Technical Support" case
# (CaseID CS2025_00073). Remove the whole row - remaining rows shift up.
$ws.Rows("24").Delete()

# Renumber CaseID column (A2:A31) sequentially as CS2025_00142 .. CS2025_00171
$lastRow = 31
for ($r = 2; $r -le $lastRow; $r++) {
    $newNum = 142 + ($r - 2)
    $newId = "CS2025_{0:D5}" -f $newNum
    $ws.Cells.Item($r, 1).Value = $newId
}
